# Fixed data: clear the stored credential values (email/password) from the
# "Book1" sheet, drop the hyperlink that pointed at the email address, and
# move the active selection to B2 - matching the sanitized fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink on A2 (mailto: link to the sample email address).
$ws.Hyperlinks.Delete()

# Clear the credential values while keeping the existing (hyperlink) style
# applied to the cells.
$ws.Range("A2").ClearContents()
$ws.Range("B2").ClearContents()

# Update the active selection to B2.
$ws.Range("B2").Select()
